# Apply crypto price/volume updates per the commit diff
# (values that look numeric are prefixed with a literal apostrophe so Excel
#  keeps storing them as text, matching the original inline-string cells)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "39.686.14"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "2.158.73"
$ws.Range("E3").Value = "  +2.52%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'227.31"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'0.628"
$ws.Range("E6").Value = "  +1.54%  "
$ws.Range("D7").Value = "'63.30"
$ws.Range("E7").Value = "  +1.59%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  +0.55%  "
$ws.Range("D10").Value = "'0.0845"
$ws.Range("E10").Value = "  +0.29%  "
$ws.Range("E11").Value = "  -0.08%  "
$ws.Range("D12").Value = "'15.91"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "2.480.62"
$ws.Range("E13").Value = "  +2.61%  "
$ws.Range("D14").Value = "'21.85"
$ws.Range("E14").Value = "  -0.96%  "
$ws.Range("E15").Value = "  -0.44%  "
$ws.Range("D16").Value = "'5.48"
$ws.Range("E16").Value = "  -0.78%  "
$ws.Range("D17").Value = "2.161.95"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "39.636.10"
$ws.Range("E18").Value = "  +1.99%  "
$ws.Range("D19").Value = "'71.67"
$ws.Range("E19").Value = "  -0.17%  "
$ws.Range("D20").Value = "'6.07"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("D21").Value = "0.0₃0843"
$ws.Range("E21").Value = "  -0.69%  "
$ws.Range("D22").Value = "'229.99"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.13%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.38"
$ws.Range("E24").Value = "  +2.76%  "
$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "'2.36"
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "'172.34"
$ws.Range("E26").Value = "  +0.39%  "
$ws.Range("D27").Value = "'9.53"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  +2.31%  "
$ws.Range("D29").Value = "'1.46"
$ws.Range("E29").Value = "  +2.48%  "
$ws.Range("D30").Value = "'19.84"
$ws.Range("E30").Value = "  +2.41%  "
$ws.Range("E31").Value = "  +4.80%  "
$ws.Range("E32").Value = "  +1.26%  "
$ws.Range("D33").Value = "'4.58"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("D34").Value = "'4.68"
$ws.Range("E34").Value = "  -1.54%  "
$ws.Range("E35").Value = "  -4.11%  "
$ws.Range("D36").Value = "'0.0618"
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("D37").Value = "'2.39"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  +2.75%  "
$ws.Range("D39").Value = "'5.13"
$ws.Range("E39").Value = "  +23.22%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").Value = "'102.64"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("E42").Value = "  -0.87%  "
$ws.Range("D43").Value = "'17.58"
$ws.Range("E43").Value = "  -2.82%  "
$ws.Range("D44").Value = "1.514.68"
$ws.Range("E44").Value = "  -0.69%  "
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'0.0919"
$ws.Range("E47").Value = "  +0.13%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").Value = "'7.75"
$ws.Range("E48").Value = "  -0.83%  "
$ws.Range("B49").Value = "ARBITRUM"
$ws.Range("C49").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D49").Value = "'1.09"
$ws.Range("E49").Value = "  -0.10%  "
$ws.Range("D50").Value = "'50.48"
$ws.Range("E50").Value = "  +9.49%  "
$ws.Range("E51").Value = "  +0.77%  "
